# Office hours note update (Week12_API slide "Announcements:")
# Replace "10-11 in person  Harsbarger 324e or by zoom on request"
# with "10-11 by zoom " while leaving the unchanged "This Thursday "
# prefix run untouched, per the class-notes edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$target = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text -like "This Thursday 10-11 in person*") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {
    # "This Thursday " (14 chars) is left alone; everything from the
    # "10-11" onward (through "...by zoom on request") is replaced.
    $prefixLen = 14
    $selection = $target.Characters($prefixLen + 1, $target.Length - $prefixLen)
    $selection.Text = "10-11 by zoom "
}
